# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 6424
$wsExpo.Range("F4").Value = 6
$wsExpo.Range("F5").Value = 381
$wsExpo.Range("F6").Value = 58
$wsExpo.Range("F9").Value = 86
$wsExpo.Range("F10").Value = 75
$wsExpo.Range("F12").Value = 156
$wsExpo.Range("F13").Value = 370
$wsExpo.Range("F14").Value = 940
$wsExpo.Range("F15").Value = 3138
$wsExpo.Range("F16").Value = 12
$wsExpo.Range("F17").Value = 190
$wsExpo.Range("F18").Value = 1807

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6424
$wsAll.Range("F4").Value = 6
$wsAll.Range("F5").Value = 381
$wsAll.Range("F6").Value = 58
$wsAll.Range("F10").Value = 86
$wsAll.Range("F11").Value = 75
$wsAll.Range("F13").Value = 156
$wsAll.Range("F14").Value = 370
$wsAll.Range("F15").Value = 940
$wsAll.Range("F16").Value = 3138
$wsAll.Range("F17").Value = 12
$wsAll.Range("F18").Value = 190
$wsAll.Range("F19").Value = 1807
